$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row labels: swap displayed names for rows 3 and 4 (A3/A4)
$ws.Range("A3").Value = "Bibi Cell Vieiralves"
$ws.Range("A4").Value = "Bibi Cell Manauara"

# Row 2 updates
$ws.Range("H2").Value = 11736.16
$ws.Range("I2").Value = 7827.31
$ws.Range("AG2").Value = 66497.78

# Row 3 updates (take previous row4 B:H values, plus new I/AG)
$ws.Range("B3").Value = 3638
$ws.Range("C3").Value = 3280.25
$ws.Range("D3").Value = 5521.8
$ws.Range("E3").Value = 2850
$ws.Range("F3").Value = 4180
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 4097.5
$ws.Range("I3").Value = 4772
$ws.Range("AG3").Value = 28339.55

# Row 4 updates (take previous row3 B:H values, plus new I/AG)
$ws.Range("B4").Value = 2756
$ws.Range("C4").Value = 3763
$ws.Range("D4").Value = 2753
$ws.Range("E4").Value = 2701
$ws.Range("F4").Value = 6020
$ws.Range("G4").Value = 3870.9
$ws.Range("H4").Value = 3677
$ws.Range("I4").Value = 1503
$ws.Range("AG4").Value = 27043.9

# Row 5 updates
$ws.Range("I5").Value = 3863.5
$ws.Range("AG5").Value = 25541.03

# Row 6 updates
$ws.Range("H6").Value = 21527.67
$ws.Range("I6").Value = 17965.81
$ws.Range("AG6").Value = 147422.26
